# Insert 3 new weekly price records (rows 175-177) for Albahaca / Vega Central
# Mapocho de Santiago, pushing the existing rows 175-195 down to 178-198.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows above the current row 175 - everything below (old
# 175..195) shifts down to 178..198, carrying its own formatting along.
$ws.Rows(175).Resize(3).Insert()

# Common (constant-across-sheet) field values for this data block.
$mercadoId = 9
$mercado   = "Vega Central Mapocho de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$catId     = 100112052
$categoria = "Albahaca"
$variedad  = "Sin especificar"
$clasif    = "Hortaliza"

function Set-Fila {
    param(
        [int]$Row,
        [double]$Fecha,
        [string]$Calidad,
        [double]$Volumen,
        [double]$PrecioMin,
        [double]$PrecioMax,
        [double]$PrecioProm,
        [string]$Unidad,
        [string]$Origen,
        [double]$PrecioKg,
        [double]$KgUnidades
    )

    $ws.Cells.Item($Row, 1).Value  = $mercadoId
    $ws.Cells.Item($Row, 2).Value  = $mercado
    $ws.Cells.Item($Row, 3).Value  = $region
    $ws.Cells.Item($Row, 4).Value  = $Fecha
    $ws.Cells.Item($Row, 5).Value  = $codreg
    $ws.Cells.Item($Row, 6).Value  = $catId
    $ws.Cells.Item($Row, 7).Value  = $categoria
    $ws.Cells.Item($Row, 8).Value  = $variedad
    $ws.Cells.Item($Row, 9).Value  = $Calidad
    $ws.Cells.Item($Row, 10).Value = $Volumen
    $ws.Cells.Item($Row, 11).Value = $PrecioMin
    $ws.Cells.Item($Row, 12).Value = $PrecioMax
    $ws.Cells.Item($Row, 13).Value = $PrecioProm
    $ws.Cells.Item($Row, 14).Value = $Unidad
    $ws.Cells.Item($Row, 15).Value = $Origen
    $ws.Cells.Item($Row, 16).Value = $PrecioKg
    $ws.Cells.Item($Row, 17).Value = $KgUnidades
    $ws.Cells.Item($Row, 18).Value = $clasif
}

# New row 175: 2021-10-07, Primera, $/docena de matas, Provincia de Chacabuco
Set-Fila 175 44476 "Primera" 61 7000 7000 7000 "$/docena de matas" "Provincia de Chacabuco" 1167 6

# New row 176: 2021-10-07, Primera, $/paquete, Región de Arica y Parinacota
Set-Fila 176 44476 "Primera" 160 4500 5000 4750 "$/paquete" "Región de Arica y Parinacota" 4750 1

# New row 177: 2021-10-07, Segunda, $/paquete, Región de Arica y Parinacota
Set-Fila 177 44476 "Segunda" 97 4000 4000 4000 "$/paquete" "Región de Arica y Parinacota" 4000 1
